$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s0 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "311.98"
$ws.Range("D2").Style = $s0
$s1 = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.69%"
$ws.Range("E2").Style = $s1
$s2 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.69"
$ws.Range("D3").Style = $s2
$s3 = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.65%"
$ws.Range("E3").Style = $s3
$s4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.129"
$ws.Range("D4").Style = $s4
$s5 = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.11%"
$ws.Range("E4").Style = $s5
$s6 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07890"
$ws.Range("D5").Style = $s6
$s7 = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.86%"
$ws.Range("E5").Style = $s7
$s8 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.414"
$ws.Range("D6").Style = $s8
$s9 = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.81%"
$ws.Range("E6").Style = $s9
$s10 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.905"
$ws.Range("D7").Style = $s10
$s11 = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.13%"
$ws.Range("E7").Style = $s11
$s12 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.275"
$ws.Range("D8").Style = $s12
$s13 = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.12%"
$ws.Range("E8").Style = $s13
$s14 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.859"
$ws.Range("D9").Style = $s14
$s15 = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-10.29%"
$ws.Range("E9").Style = $s15
$s16 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9204"
$ws.Range("D10").Style = $s16
$s17 = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.33%"
$ws.Range("E10").Style = $s17
$s18 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1180"
$ws.Range("D11").Style = $s18
$s19 = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.95%"
$ws.Range("E11").Style = $s19
$s20 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1931"
$ws.Range("D12").Style = $s20
$s21 = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.23%"
$ws.Range("E12").Style = $s21
$s22 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09037"
$ws.Range("D13").Style = $s22
$s23 = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.48%"
$ws.Range("E13").Style = $s23
$s24 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03351"
$ws.Range("D14").Style = $s24
$s25 = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.51%"
$ws.Range("E14").Style = $s25
$s26 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09606"
$ws.Range("D15").Style = $s26
$s27 = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.02%"
$ws.Range("E15").Style = $s27
$s28 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001393"
$ws.Range("D16").Style = $s28
$s29 = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.85%"
$ws.Range("E16").Style = $s29
$s30 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006037"
$ws.Range("D17").Style = $s30
$s31 = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.91%"
$ws.Range("E17").Style = $s31
$s32 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.550"
$ws.Range("D18").Style = $s32
$s33 = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.03%"
$ws.Range("E18").Style = $s33
$s34 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.292"
$ws.Range("D20").Style = $s34
$s35 = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.55%"
$ws.Range("E20").Style = $s35
$s36 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1285"
$ws.Range("D21").Style = $s36
$s37 = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.49%"
$ws.Range("E21").Style = $s37
$s38 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2586"
$ws.Range("D22").Style = $s38
$s39 = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.88%"
$ws.Range("E22").Style = $s39
$s40 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04372"
$ws.Range("D23").Style = $s40
$s41 = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.05%"
$ws.Range("E23").Style = $s41
$s42 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001247"
$ws.Range("D24").Style = $s42
$s43 = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.13%"
$ws.Range("E24").Style = $s43
$s44 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004672"
$ws.Range("D25").Style = $s44
$s45 = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "10.43%"
$ws.Range("E25").Style = $s45
$s46 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001358"
$ws.Range("D26").Style = $s46
$s47 = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.70%"
$ws.Range("E26").Style = $s47
$s48 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003985"
$ws.Range("D27").Style = $s48
$s49 = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-98.10%"
$ws.Range("E27").Style = $s49
$s50 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02260"
$ws.Range("D39").Style = $s50
$s51 = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.41%"
$ws.Range("E39").Style = $s51
$s52 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05091"
$ws.Range("D40").Style = $s52
$s53 = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.84%"
$ws.Range("E40").Style = $s53
$s54 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007448"
$ws.Range("D41").Style = $s54
$s55 = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.67%"
$ws.Range("E41").Style = $s55
$s56 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.008993"
$ws.Range("D42").Style = $s56
$s57 = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-9.72%"
$ws.Range("E42").Style = $s57
$s58 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1356"
$ws.Range("D43").Style = $s58
$s59 = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.29%"
$ws.Range("E43").Style = $s59
$s60 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001948"
$ws.Range("D44").Style = $s60
$s61 = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.38%"
$ws.Range("E44").Style = $s61
$s62 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008598"
$ws.Range("D45").Style = $s62
$s63 = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.37%"
$ws.Range("E45").Style = $s63
$s64 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006579"
$ws.Range("D46").Style = $s64
$s65 = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.85%"
$ws.Range("E46").Style = $s65
$s66 = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.04%"
$ws.Range("E47").Style = $s66
$s67 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0009989"
$ws.Range("D48").Style = $s67
$s68 = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-23.20%"
$ws.Range("E48").Style = $s68
$s69 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003057"
$ws.Range("D49").Style = $s69
$s70 = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "2.12%"
$ws.Range("E49").Style = $s70
$s71 = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("E50").Style = $s71
$s72 = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.04%"
$ws.Range("E51").Style = $s72
